$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-29 Monday" "2025-09-30 Tuesday"

Replace-Text "718×2=1436" "594×9=5346"
Replace-Text "432×5=2160" "426×6=2556"
Replace-Text "567×7=3969" "674×8=5392"
Replace-Text "915×7=6405" "396×2=792"
Replace-Text "265×7=1855" "890×3=2670"

Replace-Text "875×7=6125" "707×8=5656"
Replace-Text "213×6=1278" "202×5=1010"
Replace-Text "388×9=3492" "396×5=1980"
Replace-Text "649×6=3894" "781×8=6248"
Replace-Text "403×9=3627" "680×4=2720"

Replace-Text "215×7=1505" "799×8=6392"
Replace-Text "674×6=4044" "953×8=7624"
Replace-Text "461×8=3688" "669×4=2676"
Replace-Text "572×7=4004" "413×8=3304"
Replace-Text "409×9=3681" "845×7=5915"

Replace-Text "325×5=1625" "430×2=860"
Replace-Text "364×7=2548" "358×2=716"
Replace-Text "902×3=2706" "525×7=3675"
Replace-Text "471×6=2826" "173×4=692"
Replace-Text "852×8=6816" "802×2=1604"

Replace-Text "309×4=1236" "579×5=2895"
Replace-Text "880×8=7040" "261×8=2088"
Replace-Text "296×7=2072" "409×4=1636"
Replace-Text "665×9=5985" "495×3=1485"
Replace-Text "440×3=1320" "708×3=2124"
